$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose Date_1 column (A) moves from 2026/01/07 to 2026/01/08
$dateRows = @(2, 8, 14, 20, 26, 32, 38, 44, 50, 56, 62, 68, 74)

foreach ($r in $dateRows) {
    $cell = $ws.Range("A$r")
    $cell.NumberFormat = "@"
    $cell.Value = "2026/01/08"
    $cell.ClearFormats()
}

# Rows whose EBITDA column (B) value changes as well
$bvals = @{
    2  = "7.84"
    8  = "8.90"
    14 = "3.14"
    20 = "13.14"
    26 = "11.62"
    32 = "27.98"
    44 = "13.68"
    50 = "11.53"
    56 = "32.10"
    68 = "12.80"
    74 = "18.20"
}

foreach ($r in $bvals.Keys) {
    $cell = $ws.Range("B$r")
    $cell.NumberFormat = "@"
    $cell.Value = $bvals[$r]
    $cell.ClearFormats()
}
